$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = '[''Capitalism'', ''neoliberalism'', ''privatisation'', ''privatise'']'
$ws.Range("E2").Value = '[''capitalism,neoliberalism,privatisation,privatise'']'
$ws.Range("C3").Value = '[''Abuse'', ''against xenophobia'']'
$ws.Range("E3").Value = '[''abuse,against'', ''xenophobia'']'
$ws.Range("C4").Value = '[''poverty'', ''unemployment'', ''jobs'']'
$ws.Range("E4").Value = '[''poverty,unemployment,jobs'']'
$ws.Range("C5").Value = '[''contracts expire'', ''project end'', ''completed'']'
$ws.Range("E5").Value = '[''contract'', ''expire,project'', ''end,completed'']'
$ws.Range("C6").Value = '[''their rights'', ''recognition'']'
$ws.Range("E6").Value = '[''rights,recognition'']'
$ws.Range("C7").Value = '[''destroyed'', ''evicted'', ''demolished'', ''relocated'', ''removal'']'
$ws.Range("E7").Value = '[''destroyed,evicted,demolished,relocated,removal'']'
$ws.Range("C8").Value = '[''School'', ''university'', ''students'']'
$ws.Range("D8").Value = '[''fees'', ''costs'', ''teacher'', ''permission'', ''results'', ''policy'', ''allowed'', ''shortage'', ''not used'', ''closed'', ''residence'']'
$ws.Range("E8").Value = '[''school,university,students'']'
$ws.Range("F8").Value = '[''fees,costs,teacher,permission,results,policy,allowed,shortage,not'', ''used,closed,residence'']'
$ws.Range("D9").Value = '[''outcome'', ''result'', ''winner'', ''unfair'', ''cheat'', ''wrong'', ''councillor  '']'
$ws.Range("F9").Value = '[''outcome,result,winner,unfair,cheat,wrong,councillor'']'
$ws.Range("C10").Value = '[''Electricity'', ''Power'', ''connections'']'
$ws.Range("D10").Value = '[''cost'', ''price'', ''supply'', ''cut'', ''disconnect'']'
$ws.Range("E10").Value = '[''electricity,power,connections'']'
$ws.Range("F10").Value = '[''cost,price,supply,cut,disconnect'']'
$ws.Range("C11").Value = '[''Work'', ''workers'', ''Company'', ''Employ'', ''employer'', ''employee'', ''Labour'', ''Industrial'', ''Bosses'', ''strike'', ''management'']'
$ws.Range("D11").Value = '[''corrupt'', ''fired'', ''dismissed'', ''policy'', ''other people'', ''hours'', ''overtime'', ''equipment'', ''tools'', ''salary'', ''wages'', ''cheque'', ''food'', ''water'', ''transport'', ''housing'', ''accommodation'', ''health'', ''safety'', ''living conditions'', ''management'', ''race'', ''racism'', ''racist'', ''discrimination'', ''subsidy'', ''allowance'', ''training'', ''white'', ''promotion'', ''working conditions'', ''increase'', ''unfair'', ''poor'', ''injustice'', ''unjust'']'
$ws.Range("E11").Value = '[''work,workers,company,employ,employer,employee,labour,industrial,bosses,strike,management'']'
$ws.Range("F11").Value = '[''corrupt,fired,dismissed,policy,other'', ''people,hours,overtime,equipment,tools,salary,wages,cheque,food,water,transport,housing,accommodation,health,safety,living'', ''conditions,management,race,racism,racist,discrimination,subsidy,allowance,training,white,promotion,working'', ''conditions,increase,unfair,poor,injustice,unjust'']'
$ws.Range("C12").Value = '[''Foreigners'', ''somali'', ''Zimbabwe'', ''xenophobic'', ''xenophobia'']'
$ws.Range("E12").Value = '[''foreigners,somali,zimbabwe,xenophobic,xenophobia'']'
$ws.Range("C13").Value = '[''healthcare'', ''nurses'', ''hospital'', ''clinic'', ''doctors'']'
$ws.Range("E13").Value = '[''healthcare,nurses,hospital,clinic,doctors'']'
$ws.Range("C14").Value = '[''Housing'', ''land'', ''stands'']'
$ws.Range("E14").Value = '[''housing,land,stands'']'
$ws.Range("C15").Value = '[''Embassy'', ''solidarily with the people of'', ''war'', ''the situation in'']'
$ws.Range("E15").Value = '[''embassy,solidarily'', ''people'', ''of,war,the'', ''situation'']'
$ws.Range("D16").Value = '[''suspect'', ''criminal'']'
$ws.Range("F16").Value = '[''suspect,criminal'']'
$ws.Range("D17").Value = '[''attend'', ''solidarity'', ''gather'', ''demonstrate'', ''crowd'', ''angry'', ''support'']'
$ws.Range("F17").Value = '[''attend,solidarity,gather,demonstrate,crowd,angry,support'']'
$ws.Range("C18").Value = '[''Labour broker'', ''labour brokers'']'
$ws.Range("E18").Value = '[''labour'', ''broker,labour'', ''broker'']'
$ws.Range("C19").Value = '[''Licence'', ''permit'', ''permits'', ''licenses'']'
$ws.Range("E19").Value = '[''licence,permit,permits,licenses'']'
$ws.Range("C20").Value = '[''Mining'', ''pollution'', ''dirty'', ''air'', ''noise'', ''cracking'', ''cracked'', ''fracking '']'
$ws.Range("E20").Value = '[''mining,pollution,dirty,air,noise,cracking,cracked,fracking'']'
$ws.Range("C21").Value = '[''Not arrive'', ''not come'', ''failed to arrive'', ''failed to come'', ''not respond'', ''not answer'']'
$ws.Range("E21").Value = '[''arrive,not'', ''come,failed'', ''arrive,failed'', ''come,not'', ''respond,not'', ''answer'']'
$ws.Range("D22").Value = '[''Bills'', ''billing'', ''rates'', ''council tax'', ''corrupt'', ''corruption'', ''illegal'', ''council'', ''councillor'', ''mayor'', ''failed'', ''waited'', ''long time'', ''years'', ''consultation'', ''feedback'', ''information'', ''old councillor'', ''reinstated'', ''reinstate'', ''reappointed'', ''put back'', ''toilets'', ''sanitation'', ''pipes'', ''water'']'
$ws.Range("F22").Value = '[''bills,billing,rates,council'', ''tax,corrupt,corruption,illegal,council,councillor,mayor,failed,waited,long'', ''time,years,consultation,feedback,information,old'', ''councillor,reinstated,reinstate,reappointed,put'', ''back,toilets,sanitation,pipes,water'']'
$ws.Range("C23").Value = '[''Government'', ''National'']'
$ws.Range("E23").Value = '[''government,national'']'
$ws.Range("C24").Value = '[''Inspectors'', ''expired'']'
$ws.Range("E24").Value = '[''inspectors,expired'']'
$ws.Range("C26").Value = '[''Police'', ''SAPS'']'
$ws.Range("D26").Value = '[''crime'', ''policing'', ''incident'', ''drugs'', ''gangs'', ''arrests'', ''suspects'', ''wrongful'', ''violence'']'
$ws.Range("E26").Value = '[''police,saps'']'
$ws.Range("F26").Value = '[''crime,policing,incident,drugs,gangs,arrests,suspects,wrongful,violence'']'
$ws.Range("C27").Value = '[''Disaster'', ''flood'', ''collapse'', ''tornado'', ''storm'']'
$ws.Range("E27").Value = '[''disaster,flood,collapse,tornado,storm'']'
$ws.Range("C28").Value = '[''water'', ''poor roads'', ''quality of roads'', ''quality of the roads'', ''quality roads'']'
$ws.Range("E28").Value = '[''water,poor'', ''roads,quality'', ''roads,quality'', ''roads,quality'', ''road'']'
$ws.Range("C29").Value = '[''service delivery'']'
$ws.Range("C30").Value = '[''New law'', ''legislation'', ''proposed'']'
$ws.Range("E30").Value = '[''new'', ''law,legislation,proposed'']'
$ws.Range("C31").Value = '[''to be part of'', ''incorporated'', ''demarcation'', ''located'', ''boundary'']'
$ws.Range("E31").Value = '[''part'', ''of,incorporated,demarcation,located,boundary'']'
$ws.Range("C32").Value = '[''Premier'']'
$ws.Range("C33").Value = '[''grant'']'
$ws.Range("C34").Value = '[''Tribal court'', ''imbizo'', ''Kgotla'', ''chief'', ''kgosi'', ''induna'']'
$ws.Range("D34").Value = '[''problem'', ''dispute'', ''challenge'', ''fight'', ''protest'', ''gather'', ''demonstrate'']'
$ws.Range("E34").Value = '[''tribal'', ''court,imbizo,kgotla,chief,kgosi,induna'']'
$ws.Range("F34").Value = '[''problem,dispute,challenge,fight,protest,gather,demonstrate'']'
$ws.Range("C35").Value = '[''Ratepayers'', ''taxpayers'', ''association'']'
$ws.Range("E35").Value = '[''ratepayers,taxpayers,association'']'
$ws.Range("C36").Value = '[''Witchcraft'', ''from the dead'', ''muti'']'
$ws.Range("E36").Value = '[''witchcraft,from'', ''dead,muti'']'
$ws.Range("C37").Value = '[''Empowerment'', ''Rights'']'
$ws.Range("D37").Value = '[''female'', ''women'', ''woman'']'
$ws.Range("E37").Value = '[''empowerment,rights'']'
$ws.Range("F37").Value = '[''female,women,woman'']'
